$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest cryptos snapshot.
# Values that look like plain numbers (e.g. "258.20") are written with a
# leading quote to force text (avoids Excel silently dropping trailing
# zeros / coercing to a Number), then the style is reset to "Normal" so the
# cell keeps its original (default) formatting instead of picking up a
# text-quote-prefix style.

$ws.Range("D2").Value = "37.348.64"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.012.94"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'258.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.01%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'56.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.11%  "
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").Value = "'14.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.12%  "
$ws.Range("D13").Value = "2.309.17"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "'0.806"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("E15").Value = "  -7.30%  "
$ws.Range("E16").Value = "  -3.43%  "
$ws.Range("D17").Value = "2.016.34"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "37.269.32"
$ws.Range("D19").Value = "'69.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "0.0₃0837"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'228.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "'2.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.78%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'164.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "'9.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").Value = "'19.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'0.130"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.57%  "
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "'0.0649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "'4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").Value = "'2.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").Value = "'5.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  +3.82%  "
$ws.Range("D41").Value = "'1.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  -5.17%  "
$ws.Range("D44").Value = "1.393.89"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").Value = "'90.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "'15.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.56%  "
$ws.Range("D47").Value = "'1.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'7.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D50").Value = "2.200.80"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'1.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.88%  "
